$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.044248342514038
$ws.Range("B1").Value = 1.120322823524475
$ws.Range("D1").Value = 1.636504650115967
$ws.Range("E1").Value = 0.9877536296844482
